$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Price" column (D) values - force text format to avoid numeric auto-conversion
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.980.03'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.317.08'
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '542.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.577'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.315.73'
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.51'
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.332'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.41'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.731.31'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '59.991.19'
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.320.31'
$ws.Range("D18").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '312.98'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.58'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.76'
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.18'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '171.38'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.70'
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.85'
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.378'
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.68'
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.02'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '318.07'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '37.83'
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '136.53'
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0938'
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.79'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0491'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0222'
$ws.Range("D50").Style = "Normal"

# Update "Volume(1h)" column (E) values
$ws.Range("E2").Value = '  +2.09%  '
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("E5").Value = '  +0.78%  '
$ws.Range("E6").Value = '  -0.95%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -1.95%  '
$ws.Range("E9").Value = '  +0.51%  '
$ws.Range("E10").Value = '  +0.61%  '
$ws.Range("E11").Value = '  +0.62%  '
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("E13").Value = '  -0.41%  '
$ws.Range("E14").Value = '  -1.38%  '
$ws.Range("E15").Value = '  +0.60%  '
$ws.Range("E16").Value = '  +2.33%  '
$ws.Range("E17").Value = '  -0.63%  '
$ws.Range("E18").Value = '  +0.98%  '
$ws.Range("E19").Value = '  -0.98%  '
$ws.Range("E20").Value = '  -1.76%  '
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("E22").Value = '  -0.71%  '
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("E24").Value = '  +0.42%  '
$ws.Range("E25").Value = '  +2.03%  '
$ws.Range("E26").Value = '  -1.16%  '
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("E28").Value = '  -1.89%  '
$ws.Range("E29").Value = '  +4.01%  '
$ws.Range("E30").Value = '  +3.31%  '
$ws.Range("E31").Value = '  +0.25%  '
$ws.Range("E32").Value = '  -1.20%  '
$ws.Range("E33").Value = '  -0.84%  '
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("E35").Value = '  +4.20%  '
$ws.Range("E36").Value = '  -2.20%  '
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("E38").Value = '  -1.23%  '
$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("E40").Value = '  -0.73%  '
$ws.Range("E41").Value = '  +7.43%  '
$ws.Range("E42").Value = '  -1.43%  '
$ws.Range("E43").Value = '  +0.38%  '
$ws.Range("E44").Value = '  -3.35%  '
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("E46").Value = '  -2.17%  '
$ws.Range("E47").Value = '  +1.45%  '
$ws.Range("E48").Value = '  +2.98%  '
$ws.Range("E49").Value = '  -0.89%  '
$ws.Range("E50").Value = '  +16.47%  '
$ws.Range("E51").Value = '  +0.22%  '
